$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44274
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 14000
$ws.Range("P2").Value = 14000
$ws.Range("S2").Value = 875

# Row 3
$ws.Range("D3").Value = 44274
$ws.Range("M3").Value = 130
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 12000
$ws.Range("S3").Value = 750

# Row 4
$ws.Range("D4").Value = 44284
$ws.Range("L4").Value = 'Especial'
$ws.Range("M4").Value = 120
$ws.Range("N4").Value = 13000
$ws.Range("O4").Value = 13000
$ws.Range("P4").Value = 13000
$ws.Range("Q4").Value = '$/caja 18 kilos'
$ws.Range("S4").Value = 722
$ws.Range("T4").Value = 18

# Row 5
$ws.Range("D5").Value = 44284
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("S5").Value = 833

# Row 6
$ws.Range("D6").Value = 44284
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 12000
$ws.Range("S6").Value = 667

# Row 7
$ws.Range("D7").Value = 44267
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 13000
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 13000
$ws.Range("S7").Value = 722

# Row 8
$ws.Range("D8").Value = 44273
$ws.Range("L8").Value = 'Especial'
$ws.Range("M8").Value = 40
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("Q8").Value = '$/caja 16 kilos'
$ws.Range("S8").Value = 938
$ws.Range("T8").Value = 16

# Row 9
$ws.Range("D9").Value = 44273
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 13000
$ws.Range("O9").Value = 13000
$ws.Range("P9").Value = 13000
$ws.Range("Q9").Value = '$/caja 16 kilos'
$ws.Range("S9").Value = 812
$ws.Range("T9").Value = 16

# Row 10
$ws.Range("D10").Value = 44273
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 10000
$ws.Range("O10").Value = 10000
$ws.Range("P10").Value = 10000
$ws.Range("Q10").Value = '$/caja 16 kilos'
$ws.Range("S10").Value = 625
$ws.Range("T10").Value = 16

# Row 11
$ws.Range("D11").Value = 44309
$ws.Range("L11").Value = 'Especial'
$ws.Range("N11").Value = 20000
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 20000
$ws.Range("Q11").Value = '$/caja 18 kilos'
$ws.Range("R11").Value = 'Provincia de Melipilla'
$ws.Range("S11").Value = 1111
$ws.Range("T11").Value = 18

# Row 12
$ws.Range("D12").Value = 44309
$ws.Range("M12").Value = 60
$ws.Range("Q12").Value = '$/caja 18 kilos'
$ws.Range("R12").Value = 'Provincia de Melipilla'
$ws.Range("S12").Value = 1000
$ws.Range("T12").Value = 18

# Row 13
$ws.Range("D13").Value = 44224
$ws.Range("M13").Value = 120
$ws.Range("N13").Value = 18000
$ws.Range("O13").Value = 18000
$ws.Range("P13").Value = 18000
$ws.Range("Q13").Value = '$/caja 16 kilos'
$ws.Range("S13").Value = 1125
$ws.Range("T13").Value = 16

# Row 14
$ws.Range("D14").Value = 44298
$ws.Range("L14").Value = 'Extra (doble especial)'
$ws.Range("M14").Value = 160
$ws.Range("N14").Value = 20000
$ws.Range("O14").Value = 20000
$ws.Range("P14").Value = 20000
$ws.Range("R14").Value = 'Provincia de Melipilla'
$ws.Range("S14").Value = 1111

# Row 15
$ws.Range("D15").Value = 44222
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 18000
$ws.Range("O15").Value = 18000
$ws.Range("P15").Value = 18000
$ws.Range("Q15").Value = '$/caja 16 kilos'
$ws.Range("R15").Value = 'Provincia de Limarí'
$ws.Range("S15").Value = 1125
$ws.Range("T15").Value = 16

# Row 16
$ws.Range("D16").Value = 44295
$ws.Range("L16").Value = 'Segunda'
$ws.Range("M16").Value = 130
$ws.Range("N16").Value = 10000
$ws.Range("O16").Value = 10000
$ws.Range("P16").Value = 10000
$ws.Range("S16").Value = 556

# Row 17
$ws.Range("D17").Value = 44292
$ws.Range("M17").Value = 150
$ws.Range("N17").Value = 16000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 16000
$ws.Range("S17").Value = 889

# Row 18
$ws.Range("D18").Value = 44292
$ws.Range("M18").Value = 80
$ws.Range("N18").Value = 14000
$ws.Range("O18").Value = 14000
$ws.Range("P18").Value = 14000
$ws.Range("S18").Value = 778

# Row 19
$ws.Range("D19").Value = 44350
$ws.Range("M19").Value = 60
$ws.Range("N19").Value = 24000
$ws.Range("O19").Value = 24000
$ws.Range("P19").Value = 24000
$ws.Range("R19").Value = 'Provincia de Limarí'
$ws.Range("S19").Value = 1333

# Row 20
$ws.Range("D20").Value = 44291
$ws.Range("M20").Value = 250
$ws.Range("N20").Value = 18000
$ws.Range("O20").Value = 18000
$ws.Range("P20").Value = 18000
$ws.Range("S20").Value = 1000

# Row 21
$ws.Range("D21").Value = 44300
$ws.Range("L21").Value = 'Especial'
$ws.Range("M21").Value = 120
$ws.Range("N21").Value = 18000
$ws.Range("O21").Value = 18000
$ws.Range("P21").Value = 18000
$ws.Range("S21").Value = 1000

# Row 22
$ws.Range("D22").Value = 44300
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = 16000
$ws.Range("O22").Value = 16000
$ws.Range("P22").Value = 16000
$ws.Range("S22").Value = 889

# Row 23
$ws.Range("D23").Value = 44252
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 140
$ws.Range("N23").Value = 13000
$ws.Range("O23").Value = 13000
$ws.Range("P23").Value = 13000
$ws.Range("S23").Value = 722

# Row 24
$ws.Range("D24").Value = 44271
$ws.Range("M24").Value = 60
$ws.Range("N24").Value = 15000
$ws.Range("O24").Value = 15000
$ws.Range("P24").Value = 15000
$ws.Range("S24").Value = 833

# Row 25
$ws.Range("D25").Value = 44315
$ws.Range("M25").Value = 50
$ws.Range("N25").Value = 24000
$ws.Range("O25").Value = 24000
$ws.Range("P25").Value = 24000
$ws.Range("R25").Value = 'Provincia de Melipilla'
$ws.Range("S25").Value = 1333

# Row 26
$ws.Range("D26").Value = 44315
$ws.Range("N26").Value = 20000
$ws.Range("O26").Value = 20000
$ws.Range("P26").Value = 20000
$ws.Range("S26").Value = 1111

# Row 27
$ws.Range("D27").Value = 44279
$ws.Range("N27").Value = 14000
$ws.Range("O27").Value = 14000
$ws.Range("P27").Value = 14000
$ws.Range("S27").Value = 778

# Row 28
$ws.Range("D28").Value = 44279
$ws.Range("M28").Value = 100
$ws.Range("N28").Value = 12000
$ws.Range("O28").Value = 12000
$ws.Range("P28").Value = 12000
$ws.Range("S28").Value = 667

# Row 29
$ws.Range("D29").Value = 44330
$ws.Range("L29").Value = 'Primera'
$ws.Range("M29").Value = 50
$ws.Range("N29").Value = 23000
$ws.Range("O29").Value = 23000
$ws.Range("P29").Value = 23000
$ws.Range("S29").Value = 1278

# Row 30
$ws.Range("D30").Value = 44301

# Row 33
$ws.Range("D33").Value = 44277
$ws.Range("N33").Value = 15000
$ws.Range("O33").Value = 15000
$ws.Range("P33").Value = 15000
$ws.Range("Q33").Value = '$/caja 18 kilos'
$ws.Range("R33").Value = 'Provincia de Limarí'
$ws.Range("S33").Value = 833
$ws.Range("T33").Value = 18

# Row 34
$ws.Range("D34").Value = 44258
$ws.Range("M34").Value = 100
$ws.Range("N34").Value = 14000
$ws.Range("O34").Value = 14000
$ws.Range("P34").Value = 14000
$ws.Range("Q34").Value = '$/caja 18 kilos'
$ws.Range("R34").Value = 'Provincia de Limarí'
$ws.Range("S34").Value = 778
$ws.Range("T34").Value = 18
